$wb = $excel.ActiveWorkbook

$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# Insert a new "State" column right after "Hotel_Name" in hotel_info sheet.
$hotelSheet.Columns.Item(2).Insert()

$hotelSheet.Cells.Item(1, 2).Value = "State"
$hotelSheet.Cells.Item(2, 2).Value = "Louisiana"

# Reorder sheet tabs: review_info should come before hotel_info.
$reviewSheet.Move($hotelSheet)
